$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "59.820.79"
$ws.Range("E2").Value = "  +1.04%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.304.43"
$ws.Range("E3").Value = "  -0.89%  "

# Row 4
$ws.Range("E4").Value = "  +0.03%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "540.61"
$ws.Range("E5").Value = "  -0.26%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "129.40"
$ws.Range("E6").Value = "  -2.07%  "

# Row 7
$ws.Range("E7").Value = "  -0.04%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.571"
$ws.Range("E8").Value = "  -2.23%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.304.02"
$ws.Range("E9").Value = "  -0.72%  "

# Row 10
$ws.Range("E10").Value = "  +0.23%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.55"
$ws.Range("E11").Value = "  +0.50%  "

# Row 12
$ws.Range("E12").Value = "  -0.48%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.331"
$ws.Range("E13").Value = "  -0.64%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "23.23"
$ws.Range("E14").Value = "  -2.60%  "

# Row 15
$ws.Range("B15").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C15").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.721.73"
$ws.Range("E15").Value = "  -0.60%  "

# Row 16
$ws.Range("B16").Value = "WrappedBTC"
$ws.Range("C16").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "59.851.65"
$ws.Range("E16").Value = "  +1.08%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000131"
$ws.Range("E17").Value = "  -0.86%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.316.42"
$ws.Range("E18").Value = "  -0.24%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.44"
$ws.Range("E19").Value = "  -1.56%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.04"
$ws.Range("E20").Value = "  -3.45%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "310.07"
$ws.Range("E21").Value = "  -1.31%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.52"
$ws.Range("E22").Value = "  -1.69%  "

# Row 23
$ws.Range("E23").Value = "  -0.29%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "63.28"

# Row 25
$ws.Range("E25").Value = "  -3.35%  "

# Row 26
$ws.Range("E26").Value = "  -0.14%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.72"
$ws.Range("E27").Value = "  -3.50%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.34"
$ws.Range("E28").Value = "  +1.73%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.19"
$ws.Range("E29").Value = "  +1.27%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "171.74"
$ws.Range("E30").Value = "  +0.24%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.70"
$ws.Range("E31").Value = "  -1.10%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0₃0721"
$ws.Range("E32").Value = "  -2.78%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.81"
$ws.Range("E33").Value = "  -1.25%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.377"
$ws.Range("E34").Value = "  -2.11%  "

# Row 35
$ws.Range("E35").Value = "  -7.25%  "

# Row 36
$ws.Range("E36").Value = "  +0.01%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "17.63"
$ws.Range("E37").Value = "  -1.31%  "

# Row 38
$ws.Range("E38").Value = "  -0.20%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.99"
$ws.Range("E39").Value = "  -1.75%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "316.99"
$ws.Range("E40").Value = "  +0.01%  "

# Row 41
$ws.Range("B41").Value = "OKB"
$ws.Range("C41").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "37.42"
$ws.Range("E41").Value = "  -1.63%  "

# Row 42
$ws.Range("B42").Value = "Stacks"
$ws.Range("C42").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.51"
$ws.Range("E42").Value = "  -1.35%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "136.23"
$ws.Range("E43").Value = "  -4.22%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.42"
$ws.Range("E44").Value = "  -0.65%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0943"
$ws.Range("E45").Value = "  -1.18%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.567"
$ws.Range("E46").Value = "  +1.39%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "18.61"
$ws.Range("E47").Value = "  +1.35%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0488"
$ws.Range("E48").Value = "  -1.18%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0₆0221"
$ws.Range("E49").Value = "  +21.11%  "

# Row 50
$ws.Range("E50").Value = "  +0.97%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "11.00"
$ws.Range("E51").Value = "  -0.08%  "
